# Added two new Mac-Addresses (10 new device rows, ids 3000166-3000175)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids         = @(3000166, 3000167, 3000168, 3000169, 3000170, 3000171, 3000172, 3000173, 3000174, 3000175)
$names       = @("Finger Print Scanner 30", "IRIS Scanner 30", "Web Camera 30", "Document Scanner 30", "Printer 30", "Finger Print Scanner 31", "IRIS Scanner 31", "Web Camera 31", "Document Scanner 31", "Printer 31")
$macs        = @("D6-15-AC-80-6B-86", "6D-58-E2-DF-74-34", "E2-A8-56-86-15-30", "72-E8-B9-FD-63-65", "D3-F3-A4-50-AD-12", "06-16-D0-0B-A6-E4", "21-78-45-AC-E9-20", "3C-E8-87-99-DB-FA", "BF-55-53-98-40-08", "5A-43-36-46-22-EB")
$serials     = @("BS563Q2230814", "BS563Q2230815", "BS563Q2230816", "BS563Q2230817", "BS563Q2230818", "BS563Q2230819", "BS563Q2230820", "BS563Q2230821", "BS563Q2230822", "BS563Q2230823")
$dspecIds    = @(165, 327, 736, 801, 920, 165, 327, 736, 801, 920)

$startRow = 147
for ($i = 0; $i -lt $ids.Length; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $ids[$i]
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $macs[$i]
    $ws.Cells.Item($r, 4).Value = $serials[$i]
    $ws.Cells.Item($r, 6).Value = $dspecIds[$i]
    $ws.Cells.Item($r, 7).Value = "eng"

    $ws.Cells.Item($r, 8).Value = $true
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131

    $ws.Cells.Item($r, 9).Value = "superadmin"
    $ws.Cells.Item($r, 10).Value = "now()"
    $ws.Cells.Item($r, 11).Value = "now()"
}

[void]$ws.Range("D145").Select()

Write-Output "Added rows 147-156"
